# "FIXED :: Excel File / BaseBall is now more less!"
# Rebalance the sword-forge chance table: the BaseBall cost curve (column J,
# and derived K/L) is scaled way down, the N7/O7 "bonus" cells become
# formulas driven by new ASBlevel/RBlevel inputs at N11/O11 (labeled via
# N10/O10), and the view is re-pointed/re-zoomed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column J (BaseBall) rebalance ------------------------------------
# Rows 4-23 and 25-43: J = 200 + 10*(row-4)  (was 1000 + 50*(row-4))
# Row 24 is a text divider row ("-") and is left untouched.
for ($r = 4; $r -le 23; $r++) {
    $ws.Cells.Item($r, 10).Value = 200 + 10*($r - 4)
}
for ($r = 25; $r -le 43; $r++) {
    $ws.Cells.Item($r, 10).Value = 200 + 10*($r - 4)
}

# Rows 45-54: J = 600 + 20*(row-44)  (was 3000 + 100*(row-45))
# Row 44 is a text divider row ("-") and is left untouched.
for ($r = 45; $r -le 54; $r++) {
    $ws.Cells.Item($r, 10).Value = 600 + 20*($r - 44)
}

# --- N7 / O7 become formulas off the new level inputs ------------------
$ws.Range("N7").Formula = '=$N$11*32'
$ws.Range("O7").Formula = '=$O$11*16'

# --- New ASBlevel / RBlevel labeled inputs at N10:O11 -------------------
$ws.Range("N10").Value = "ASBlevel"
$ws.Range("O10").Value = "RBlevel"
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 5

# --- N4 halves alongside the rest of the curve --------------------------
$ws.Range("N4").Value = 200

# --- View: re-zoomed and re-selected ------------------------------------
$ws.Activate()
$ws.Range("N12").Select()
$excel.ActiveWindow.Zoom = 57
